$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 49600000
$ws.Range("I8").Value = 49600000
$ws.Range("K8").Value = 148800000
$ws.Range("M8").Value = -148799861

$ws.Range("H40").Value = 3954
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -825

$ws.Range("H51").Value = 41669116
$ws.Range("J51").Value = 41669116
$ws.Range("L51").Value = 41669116
$ws.Range("N51").Value = -41670084

$ws.Range("H80").Value = 725.6923
$ws.Range("I80").Value = 602.1667
$ws.Range("J80").Value = 831.5714
$ws.Range("K80").Value = 1806.5001
$ws.Range("L80").Value = 2494.7142
$ws.Range("M80").Value = -808.5001
$ws.Range("N80").Value = -4490.7142

$ws.Range("H83").Value = 725.6923
$ws.Range("I83").Value = 602.1667
$ws.Range("J83").Value = 831.5714
$ws.Range("K83").Value = 5419.5003
$ws.Range("L83").Value = 7484.1426
$ws.Range("M83").Value = -427.5002999999997
$ws.Range("N83").Value = -17468.1426

$ws.Range("H96").Value = 1350.8235
$ws.Range("I96").Value = 701.5
$ws.Range("J96").Value = 1705
$ws.Range("K96").Value = 2104.5
$ws.Range("L96").Value = 5115
$ws.Range("M96").Value = -731.5
$ws.Range("N96").Value = -7861

$ws.Range("H137").Value = 4468.096
$ws.Range("I137").Value = 3675.2354
$ws.Range("J137").Value = 5965.722
$ws.Range("K137").Value = 11025.7062
$ws.Range("L137").Value = 17897.166
$ws.Range("M137").Value = -8475.706200000001
$ws.Range("N137").Value = -22997.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 362723.6
$ws.Range("I45").Value = 482574.16
$ws.Range("K45").Value = 482574.16
$ws.Range("M45").Value = -482197.16

$ws.Range("H61").Value = 1060.3334
$ws.Range("I61").Value = 954.2308
$ws.Range("K61").Value = 954.2308
$ws.Range("M61").Value = -742.2308

$ws.Range("H74").Value = 1994.5454
$ws.Range("I74").Value = 1463.3334
$ws.Range("K74").Value = 1463.3334
$ws.Range("M74").Value = -589.3334

$ws.Range("H77").Value = 1994.5454
$ws.Range("I77").Value = 1463.3334
$ws.Range("K77").Value = 7316.666999999999
$ws.Range("M77").Value = -2948.666999999999

$ws.Range("H136").Value = 1060.3334
$ws.Range("I136").Value = 954.2308
$ws.Range("K136").Value = 2862.6924
$ws.Range("M136").Value = -312.6923999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1138.3636
$ws.Range("I134").Value = 1138.3636
$ws.Range("K134").Value = 3415.0908
$ws.Range("M134").Value = -880.0907999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 145128.28
$ws.Range("J31").Value = 2932.6667
$ws.Range("L31").Value = 2932.6667
$ws.Range("N31").Value = -3522.6667

$ws.Range("H34").Value = 145128.28
$ws.Range("J34").Value = 2932.6667
$ws.Range("L34").Value = 2932.6667
$ws.Range("N34").Value = -3336.6667

$ws.Range("H58").Value = 2324.4375
$ws.Range("I58").Value = 1986.2858
$ws.Range("K58").Value = 1986.2858
$ws.Range("M58").Value = -1783.2858

$ws.Range("H64").Value = 25271
$ws.Range("J64").Value = 25271
$ws.Range("L64").Value = 25271
$ws.Range("N64").Value = -25767

$ws.Range("H67").Value = 25271
$ws.Range("J67").Value = 25271
$ws.Range("L67").Value = 25271
$ws.Range("N67").Value = -26987

$ws.Range("H96").Value = 11155.75
$ws.Range("J96").Value = 11155.75
$ws.Range("L96").Value = 11155.75
$ws.Range("N96").Value = -16647.75

$ws.Range("H99").Value = 15522.5
$ws.Range("I99").Value = 16079.533
$ws.Range("K99").Value = 16079.533
$ws.Range("M99").Value = -14581.533

$ws.Range("H126").Value = 15522.5
$ws.Range("I126").Value = 16079.533
$ws.Range("K126").Value = 48238.599
$ws.Range("M126").Value = -45768.599

$ws.Range("H132").Value = 4107.1
$ws.Range("J132").Value = 4608.375
$ws.Range("L132").Value = 13825.125
$ws.Range("N132").Value = -18885.125

$ws.Range("H134").Value = 2594.0256
$ws.Range("I134").Value = 2687.606
$ws.Range("K134").Value = 8062.818000000001
$ws.Range("M134").Value = -5527.818000000001

$ws.Range("H136").Value = 2324.4375
$ws.Range("I136").Value = 1986.2858
$ws.Range("K136").Value = 5958.857400000001
$ws.Range("M136").Value = -3408.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 458.875
$ws.Range("I122").Value = 496
$ws.Range("K122").Value = 4464
$ws.Range("M122").Value = -2014

$ws.Range("H131").Value = 131837.9
$ws.Range("I131").Value = 500499.5
$ws.Range("J131").Value = 39672.5
$ws.Range("K131").Value = 1501498.5
$ws.Range("L131").Value = 119017.5
$ws.Range("M131").Value = -1496458.5
$ws.Range("N131").Value = -129097.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744

$ws.Range("H102").Value = 1426.6
$ws.Range("I102").Value = 1413.6471
$ws.Range("K102").Value = 1413.6471
$ws.Range("M102").Value = 208.3529000000001

$ws.Range("H132").Value = 1464.3667
$ws.Range("I132").Value = 1528.5769
$ws.Range("J132").Value = 1047
$ws.Range("K132").Value = 4585.7307
$ws.Range("L132").Value = 3141
$ws.Range("M132").Value = -2055.7307
$ws.Range("N132").Value = -8201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 302.4
$ws.Range("I9").Value = 128
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 128
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 96
$ws.Range("N9").Value = -1448

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H98").Value = 30355
$ws.Range("J98").Value = 30355
$ws.Range("L98").Value = 30355
$ws.Range("N98").Value = -36345

$ws.Range("H132").Value = 2918.6223
$ws.Range("I132").Value = 2441.8
$ws.Range("K132").Value = 7325.400000000001
$ws.Range("M132").Value = -4795.400000000001

$ws.Range("H136").Value = 3391.1892
$ws.Range("I136").Value = 2965
$ws.Range("K136").Value = 8895
$ws.Range("M136").Value = -6345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336
$ws.Range("M14").ClearContents()

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H122").Value = 1628.5714
$ws.Range("I122").Value = 1487.4117
$ws.Range("K122").Value = 4462.2351
$ws.Range("M122").Value = -2012.2351

$ws.Range("H132").Value = 1365273.8
$ws.Range("I132").Value = 6046.24
$ws.Range("K132").Value = 18138.72
$ws.Range("M132").Value = -15608.72

$ws.Range("H136").Value = 1671
$ws.Range("I136").Value = 910.6667
$ws.Range("K136").Value = 2732.0001
$ws.Range("M136").Value = -182.0001000000002
